# Add "NA" values under the duplicate_image_filename column (column E)
# for the data rows (rows 2-21), per commit message:
# "add the NA's under duplicate_image_filename"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2:E21").Value = "NA"
